$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2026
$ws.Range("L3").Value = 2060
$ws.Range("J4").Value = 574
$ws.Range("L4").Value = 572
$ws.Range("L6").Value = 1843
$ws.Range("J7").Value = 8262
$ws.Range("L7").Value = 6618

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 117
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 416

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 50
$ws.Range("L3").Value = 64
$ws.Range("L7").Value = 156

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 100
$ws.Range("L6").Value = 101
$ws.Range("L7").Value = 294

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 98

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 67
$ws.Range("L7").Value = 237

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L2").Value = 48
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 129

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 50
$ws.Range("L7").Value = 222
$ws.Range("L8").Value = 416
$ws.Range("L9").Value = 40
$ws.Range("L10").Value = 45
$ws.Range("L18").Value = 48
$ws.Range("L19").Value = 187
$ws.Range("L29").Value = 340
$ws.Range("L33").Value = 294
$ws.Range("L34").Value = 42
$ws.Range("L37").Value = 237
$ws.Range("L42").Value = 206
$ws.Range("L47").Value = 45
$ws.Range("L51").Value = 78
$ws.Range("L52").Value = 135
$ws.Range("L53").Value = 79
$ws.Range("L54").Value = 139
$ws.Range("J63").Value = 63
$ws.Range("L63").Value = 26
$ws.Range("L64").Value = 49
$ws.Range("L65").Value = 129
$ws.Range("L67").Value = 236
$ws.Range("L76").Value = 69
$ws.Range("L78").Value = 93
$ws.Range("L83").Value = 156
$ws.Range("L84").Value = 69
$ws.Range("L85").Value = 349
$ws.Range("L86").Value = 51
$ws.Range("L87").Value = 19
$ws.Range("L89").Value = 85
$ws.Range("L93").Value = 36
$ws.Range("L94").Value = 79
$ws.Range("L95").Value = 98
$ws.Range("L96").Value = 60
$ws.Range("L98").Value = 49
$ws.Range("L99").Value = 102
$ws.Range("J101").Value = 8262
$ws.Range("L101").Value = 6618

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 139

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 122
$ws.Range("L6").Value = 90
$ws.Range("L7").Value = 340

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 58
$ws.Range("L6").Value = 59
$ws.Range("L7").Value = 187

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 69

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L4").Value = 12
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L2").Value = 26
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 18
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 70
$ws.Range("L7").Value = 222

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 20
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("L3").Value = 9
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 49

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L2").Value = 10
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 14
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L2").Value = 9
$ws.Range("L6").Value = 51

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 78

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 108
$ws.Range("L7").Value = 349

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 41
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L2").Value = 6
$ws.Range("L7").Value = 19
